$wb = $excel.ActiveWorkbook

# --- Model sheet: change the selected livestock type from Pig to Cattle ---
$model = $wb.Worksheets.Item("Model")
$model.Range("C6").Value = "Cattle"

# --- Lookup sheet: fix MATCH() calls to use exact-match mode (0) ---
$lookup = $wb.Worksheets.Item("Lookup")
$lookup.Range("C2").Formula = "=MATCH(B2,Parameters!A3:A5,0)"
$lookup.Range("C3").Formula = "=MATCH(B3,Parameters!B2:C2,0)"

# --- ChangeLog sheet: append a new version row documenting this fix ---
$changelog = $wb.Worksheets.Item("ChangeLog")
$changelog.Range("A9").Value = 1.1
$changelog.Range("B9").NumberFormat = $changelog.Range("B8").NumberFormat
$changelog.Range("B9").Value = "4/17/2023"
$changelog.Range("C9").Value = "AMOSTO.xlsx"
$changelog.Range("D9").Value = "Sasha"
$changelog.Range("E9").Value = 'Fix lookup behavior in "Lookup" for cattle/digestate based on user feedback'

# --- Restore cursor/selection positions on the touched sheets, as left by the author ---
$parameters = $wb.Worksheets.Item("Parameters")
$parameters.Activate() | Out-Null
$parameters.Range("A5").Select() | Out-Null

$lookup.Activate() | Out-Null
$lookup.Range("C3").Select() | Out-Null

$changelog.Activate() | Out-Null
$changelog.Range("A9").Select() | Out-Null

$model.Activate() | Out-Null

$wb.Save()
